$wb = $excel.ActiveWorkbook

# "Generate Report for handback" - refresh the handoff/handback timestamps
# (column D = Correspond Handoff Datetime, column G = Correspond Handback
# DateTime) for the most recently handed-back file (row 2) on each
# per-language worksheet.

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-13 11:12:41"
$zhcn.Range("G2").Value = "2016-01-13 11:13:32"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-13 11:12:53"
$dede.Range("G2").Value = "2016-01-13 11:13:53"
